# Corregí los psp's de la ClaseConfiguración
#
# - Rename the worksheet "excel(1)" -> "excel" (this also updates the
#   defined name's RefersTo range target automatically).
# - Rename the defined name "excel_1" -> "excel".
# - Refresh the report timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet from "excel(1)" to "excel"
$ws.Name = "excel"

# Rename the defined name from "excel_1" to "excel"
# ($n.Name comes back scoped as "<sheet>!<name>", so compare on the part
# after the last "!" rather than the raw string.)
foreach ($n in $wb.Names) {
    $parts = $n.Name -split "!"
    $localName = $parts[$parts.Length - 1]
    if ($localName -eq "excel_1") {
        $n.Name = "excel"
    }
}

# Update the generated-report timestamp text
$ws.Range("A23").Value = "Reporte generado a las 11:23 AM el 5/12/2018"
